$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the closing italic paragraph's text (the old meta-description
#    sentence) with the new image-prompt text. Do this FIRST, while the old
#    sentence is still unique in the document (it gets reused below too).
# ---------------------------------------------------------------------------

$oldTail = "Read our review of the 2 Gods Zeus versus Thor online slot game and play for free. Enjoy dual spin, win spins feature and immersive theme combining Greek and Norse mythology."
$newTail = 'Create a feature image for "2 Gods Zeus versus Thor" that captures the game''s unique blend of Greek and Norse mythology, as well as its innovative Dual Spin mechanism. The image should be in a cartoon style, with bright, vivid colors that will grab players'' attention. It should feature a happy Maya warrior with glasses, symbolizing the fun and exciting gameplay of the slot game. The image should show Zeus and Thor, each on their own side of the game grid, facing off against each other in a fierce battle. The background should be a mix of Greek and Norse imagery, including lightning bolts, thunderclouds, Viking ships, and Greek temples. The Dual Spin mechanism should be prominently displayed, perhaps through the use of two different colored arrows or spin buttons. The Maya warrior should be shown standing in front of the game grid, looking excited and happy as he prepares to enter the world of mythical gods and legendary battles. He should be wearing glasses to symbolize that this is a game of strategy and skill, not just luck. Overall, the feature image should be fun, engaging, and dynamic, capturing the essence of "2 Gods Zeus versus Thor" and encouraging players to try out this exciting and innovative slot game.'

$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph ("Play 2 Gods Zeus versus Thor Slot for Free - Review").
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)

# Make a new (initially empty) paragraph right after the title.
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Seed it with a copy of the bold "Play 2 Gods..." run further down the
# document so we inherit the same leading empty run + bold-run structure,
# then rewrite the text/formatting in place.
$paraCount = $d.Paragraphs.Count
$srcIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if (($i -ne 1) -and ($ptext -eq "Play 2 Gods Zeus versus Thor Slot for Free " + [char]0x2013 + " Review" + [char]0x0D)) {
        $srcIndex = $i
    }
}
$srcPara = $d.Paragraphs($srcIndex)
$metaPara.Range.FormattedText = $srcPara.Range.FormattedText

$metaStart = $metaPara.Range.Start
$metaPara.Range.Find.Execute(
    "Play 2 Gods Zeus versus Thor Slot for Free " + [char]0x2013 + " Review",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description", 2) | Out-Null

$boldEnd = $metaStart + 16
$boldRange = $d.Range($metaStart, $boldEnd)
$boldRange.InsertAfter(": Read our review of the 2 Gods Zeus versus Thor online slot game and play for free. Enjoy dual spin, win spins feature and immersive theme combining Greek and Norse mythology.")

# ---------------------------------------------------------------------------
# 3) Remove the duplicate bold "Play 2 Gods..." paragraph that used to sit
#    right before the closing italic paragraph.
# ---------------------------------------------------------------------------

$paraCount = $d.Paragraphs.Count
$dupIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if (($i -ne 2) -and ($ptext -eq "Play 2 Gods Zeus versus Thor Slot for Free " + [char]0x2013 + " Review" + [char]0x0D)) {
        $dupIndex = $i
    }
}
if ($dupIndex -gt 0) {
    $d.Paragraphs($dupIndex).Range.Delete() | Out-Null
}

Write-Output "done"
